$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 72 (shifts existing rows 72..145 down to 73..146)
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new weekly data point.
$ws.Cells.Item(72, 1).Value = 7
$ws.Cells.Item(72, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value = "Ñuble"
$ws.Cells.Item(72, 4).Value = 44512
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 100112017
$ws.Cells.Item(72, 7).Value = "Apio"
$ws.Cells.Item(72, 8).Value = "Americana (o)"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 100
$ws.Cells.Item(72, 11).Value = 8000
$ws.Cells.Item(72, 12).Value = 9000
$ws.Cells.Item(72, 13).Value = 8500
$ws.Cells.Item(72, 14).Value = "$/docena de matas"
$ws.Cells.Item(72, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(72, 16).Value = 1417
$ws.Cells.Item(72, 17).Value = 6
$ws.Cells.Item(72, 18).Value = "Hortaliza"
